# Mark the two remaining "amazonTest" rows as executed (B8, B9: "no" -> "yes")
# and leave the selection on the last-edited cell (E7), matching the
# recorded Excel session state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$ws.Range("B8").Value = "yes"
$ws.Range("B9").Value = "yes"

$ws.Activate()
$ws.Range("E7").Select()
